$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 63 data, continuing the pattern of rows above
$ws.Range("A63").Value = 20220308
$ws.Range("B63").Value = 2220.8926944598902
$ws.Range("C63").Value = 2224.4699999999998
$ws.Range("D63").Formula = "=100*(B63-C63)/C63"
$ws.Range("E63").Value = 180
$ws.Range("F63").Value = "CRM OPENED 20220302"

# Update the active selection to match the new last entry row
$ws.Range("F60").Select()

$wb.Save()
